# Project Completed, Final Changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "ITSPOC"
$ws.Range("B1").Value = "TASKDESCRIPTION"
$ws.Range("C1").Value = "CURRENTSTAGE"
$ws.Range("D1").Value = "ACTIONPOINT"
$ws.Range("E1").Value = "RESPONSIBILITY"
$ws.Range("F1").Value = "TARGET"
$ws.Range("G1").Value = "EMAIL"

# ---- Row 2 ----
$ws.Range("A2").Value = "USER1"
$ws.Range("B2").Value = "NEW SOMETHING"
$ws.Range("C2").Value = "EXAMPLE 1"
$ws.Range("D2").Value = "EXAMPLE 1"
$ws.Range("E2").Value = "RAM"
$ws.Range("F2").Value = 44177
$ws.Range("G2").Value = "swasti.tiwari@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:swasti.tiwari@gmail.com")

# ---- Row 3 ----
$ws.Range("A3").Value = "USER2"
$ws.Range("B3").Value = "NEW SOMETHING 2"
$ws.Range("C3").Value = "example 2"
$ws.Range("D3").Value = "exampl2 "
$ws.Range("E3").Value = "abc"
$ws.Range("F3").Value = 44177
$ws.Range("G3").Value = "3as1827000176@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G3"), "mailto:3as1827000176@gmail.com")

# ---- Column widths ----
$ws.Columns("A:B").ColumnWidth = 16
$ws.Columns("G").ColumnWidth = 23.833333333333332

# ---- Selection ----
[void]$ws.Range("F3").Select()
